$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format numeric/percentage columns (D, E) as Text so values round-trip
# exactly as strings (matching the source data which stores these as text).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "317.80"
$ws.Range("E2").Value = "-3.15%"
$ws.Range("D3").Value = "42.09"
$ws.Range("E3").Value = "-4.54%"
$ws.Range("D4").Value = "5.200"
$ws.Range("E4").Value = "-3.63%"
$ws.Range("D5").Value = "0.08127"
$ws.Range("E5").Value = "-2.93%"
$ws.Range("D6").Value = "4.372"
$ws.Range("E6").Value = "-1.31%"
$ws.Range("D7").Value = "1.752"
$ws.Range("E7").Value = "-9.82%"
$ws.Range("D8").Value = "0.9320"
$ws.Range("E8").Value = "-4.32%"
$ws.Range("E9").Value = "-0.99%"
$ws.Range("D10").Value = "0.1863"
$ws.Range("E10").Value = "-1.85%"
$ws.Range("D11").Value = "0.09277"
$ws.Range("E11").Value = "-4.47%"
$ws.Range("D12").Value = "0.04579"
$ws.Range("E12").Value = "-2.26%"
$ws.Range("D13").Value = "7.407"
$ws.Range("E13").Value = "-19.03%"
$ws.Range("D14").Value = "0.1056"
$ws.Range("E14").Value = "-0.74%"
$ws.Range("D15").Value = "0.001296"
$ws.Range("E15").Value = "0.31%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.005886"
$ws.Range("E16").Value = "-4.41%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.352"
$ws.Range("E17").Value = "-1.16%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "2.545"
$ws.Range("E18").Value = "0.95%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "0.3376"
$ws.Range("E19").Value = "1.39%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "0.1382"
$ws.Range("E20").Value = "0.70%"
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D21").Value = "0.2597"
$ws.Range("E21").Value = "1.84%"
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D22").Value = "0.04178"
$ws.Range("E22").Value = "0.55%"
$ws.Range("D23").Value = "0.001244"
$ws.Range("E23").Value = "-4.10%"
$ws.Range("D24").Value = "0.004257"
$ws.Range("E24").Value = "-3.28%"
$ws.Range("D25").Value = "0.0001222"
$ws.Range("E25").Value = "-6.05%"
$ws.Range("D26").Value = "0.0002984"
$ws.Range("E26").Value = "-0.09%"
$ws.Range("D38").Value = "0.02601"
$ws.Range("E38").Value = "-2.10%"
$ws.Range("D39").Value = "0.05499"
$ws.Range("E39").Value = "-2.30%"
$ws.Range("D40").Value = "0.008056"
$ws.Range("E40").Value = "2.83%"
$ws.Range("D41").Value = "0.1394"
$ws.Range("E41").Value = "-1.25%"
$ws.Range("D42").Value = "0.006527"
$ws.Range("E42").Value = "-11.48%"
$ws.Range("D43").Value = "0.002091"
$ws.Range("E43").Value = "-0.64%"
$ws.Range("D44").Value = "0.008227"
$ws.Range("E44").Value = "-4.72%"
$ws.Range("D45").Value = "0.3473"
$ws.Range("E45").Value = "-1.10%"
$ws.Range("D46").Value = "0.00006730"
$ws.Range("E46").Value = "-1.74%"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "-0.01%"
$ws.Range("D48").Value = "0.003391"
$ws.Range("E48").Value = "-3.41%"
$ws.Range("D49").Value = "0.004107"
$ws.Range("E49").Value = "16.15%"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").Value = "-0.01%"
$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").Value = "-0.01%"
